# Fruta / hortaliza, semanal
# A new weekly price record is inserted before the current row 541 ("Feria
# Lagunitas de Puerto Montt" / "Coliflor" data block), shifting every
# subsequent row down by one (541-602 -> 542-603) and extending the used
# range from A1:R602 to A1:R603.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at position 541; this pushes the old rows
# 541..602 down to 542..603 and keeps the existing formatting (date style
# on column D) because Excel copies the format of the row above on insert.
$ws.Rows("541:541").Insert()

# Populate the newly inserted row 541 with the new weekly record.
$ws.Range("A541").Value = 4
$ws.Range("B541").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C541").Value = "Los Lagos"
$ws.Range("D541").Value = 45194
$ws.Range("E541").Value = 10
$ws.Range("F541").Value = 100112008
$ws.Range("G541").Value = "Coliflor"
$ws.Range("H541").Value = "Sin especificar"
$ws.Range("I541").Value = "Primera"
$ws.Range("J541").Value = 500
$ws.Range("K541").Value = 1500
$ws.Range("L541").Value = 1500
$ws.Range("M541").Value = 1500
$ws.Range("N541").Value = "`$/unidad"
$ws.Range("O541").Value = "Región Metropolitana"
$ws.Range("P541").Value = 1500
$ws.Range("Q541").Value = 1
$ws.Range("R541").Value = "Hortaliza"
